$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "success" header in D1, styled like the other headers (B1/C1) ---
$ws.Cells.Item(1, 4).Value = "success"
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)   # xlPasteFormats - copy header style (bold/border/centered)

# success = "1" for the rows whose list starts with ['00','01','02'] (A2:A8),
# success = "0" for the remaining rows (A9:A13)
$successValues = @("1", "1", "1", "1", "1", "1", "1", "0", "0", "0", "0", "0")

for ($i = 0; $i -lt $successValues.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 4)

    # Force the value to be stored as text (shared string "1"/"0"), matching
    # the data rows, instead of Excel's default numeric auto-detection.
    $cell.NumberFormat = "@"
    $cell.Value = $successValues[$i]

    # Reset formatting back to the plain (unstyled) look used by columns B/C
    # on data rows, same as the other plain cells in this sheet.
    $ws.Cells.Item($row, 3).Copy()
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
